$wb = $excel.ActiveWorkbook

# --- Sheet "Alunos": update the Categoria of Jacir Novais (row 21) ---
$wsAlunos = $wb.Worksheets.Item("Alunos")
$wsAlunos.Range("K21").Value = "Não definida"

# --- Sheet "Registros": fill in attendance marks for rows 6-8 ---
$wsRegistros = $wb.Worksheets.Item("Registros")

# joão do pão (row 6)
$wsRegistros.Range("AC6").Value = "c"
$wsRegistros.Range("AG6").Value = "f"
$wsRegistros.Range("AH6").Value = "f"

# fernando lando (row 7)
$wsRegistros.Range("AC7").Value = "c"
$wsRegistros.Range("AG7").Value = "c"
$wsRegistros.Range("AH7").Value = "c"

# ana cintra (row 8)
$wsRegistros.Range("AC8").Value = "c"
$wsRegistros.Range("AG8").Value = "c"
$wsRegistros.Range("AH8").Value = "j"
